$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    5  = -9
    6  = -2
    9  = -1
    12 = 0
    19 = -4
    24 = -3
    25 = 3
    27 = 0
    29 = -2
    36 = -1
    38 = 3
    45 = -1
    47 = 1
    48 = -3
    50 = 1
    53 = -1
    54 = 2
    56 = 2
    59 = 1
    68 = -6
    73 = -1
    74 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
